$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells K1/L1 -------------------------------------------------
# Copy the formatting from the existing "orig-fxp-drop" header cell (J1, style
# index 1: bold font, thin border, centered/top aligned) onto the two new
# header cells so they look consistent with the rest of the header row.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"

# --- New data columns K2:L21 ------------------------------------------------
$values = @{
    2  = @(0.8835978835978836, 0.005291005291005235)
    3  = @(0.4603174603174603, 0)
    4  = @(0.4603174603174603, 0)
    5  = @(0.4603174603174603, 0)
    6  = @(0.7513227513227513, 0.126984126984127)
    7  = @(0.8835978835978836, 0)
    8  = @(0.7513227513227513, 0.1164021164021164)
    9  = @(0.8835978835978836, -0.01058201058201058)
    10 = @(0.8835978835978836, -0.01587301587301593)
    11 = @(0.8783068783068783, 0)
    12 = @(0.8835978835978836, 0)
    13 = @(0.4603174603174603, 0)
    14 = @(0.4603174603174603, 0)
    15 = @(0.8888888888888888, 0)
    16 = @(0.8835978835978836, -0.01587301587301593)
    17 = @(0.8835978835978836, 0)
    18 = @(0.8888888888888888, -0.02116402116402116)
    19 = @(0.7936507936507936, 0.09523809523809523)
    20 = @(0.4603174603174603, 0)
    21 = @(0.7142857142857143, 0.1746031746031745)
}

foreach ($row in $values.Keys | Sort-Object) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 11).Value = $pair[0]
    $ws.Cells.Item($row, 12).Value = $pair[1]
}

Write-Output "K1:L21 populated"
